$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "I learned how to calculate the Euclidean distance between test and training data using numpy.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Calculating the Euclidean distance between test and training data using numpy.",
    2)

$d.Content.Find.Execute(
    "I learned how to find the K-nearest neighbors and decide the majority class using numpy and Counter.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Finding the K-nearest neighbors and deciding the majority class using numpy and Counter.",
    2)

$d.Content.Find.Execute(
    "I learned how to compare the prediction with the ground truth in the test data using numpy.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Comparing the prediction with the ground truth in the test data using numpy.",
    2)

$d.Content.Find.Execute(
    "I also learned how to create a Word document using the python-docx library and add tables and images to it.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Creating a Word document using the python-docx library and adding tables and images to it.(This was actually pretty cool!)",
    2)
